$wb = $excel.ActiveWorkbook

# =========================================================
# Sheet: "Translations - Common" (new Czech translation keys
# for the new "create cell" feature)
# =========================================================
$wsCommon = $wb.Worksheets.Item("Translations - Common")

# Template row used to clone the existing cell style (wrapped text, s="1")
$commonTemplate = $wsCommon.Range("A1566:C1566")

# --- rows 1567-1571 ---
$commonTemplate.Copy($wsCommon.Range("A1567:C1567"))
$wsCommon.Cells.Item(1567, 1).Value = "cs"
$wsCommon.Cells.Item(1567, 2).Value = "lab.build.atomizer.coil.favorite.tab"
$wsCommon.Cells.Item(1567, 3).Value = "Oblíbené odporové dráty"

$commonTemplate.Copy($wsCommon.Range("A1568:C1568"))
$wsCommon.Cells.Item(1568, 1).Value = "cs"
$wsCommon.Cells.Item(1568, 2).Value = "lab.cell.create.button"
$wsCommon.Cells.Item(1568, 3).Value = "Nový článek"

$commonTemplate.Copy($wsCommon.Range("A1569:C1569"))
$wsCommon.Cells.Item(1569, 1).Value = "cs"
$wsCommon.Cells.Item(1569, 2).Value = "lab.cell.create.title"
$wsCommon.Cells.Item(1569, 3).Value = "Nový článek"

$commonTemplate.Copy($wsCommon.Range("A1570:C1570"))
$wsCommon.Cells.Item(1570, 1).Value = "cs"
$wsCommon.Cells.Item(1570, 2).Value = "shared.cell.create.403.title"
$wsCommon.Cells.Item(1570, 3).Value = "K této funkci nemáte přístup."

$commonTemplate.Copy($wsCommon.Range("A1571:C1571"))
$wsCommon.Cells.Item(1571, 1).Value = "cs"
$wsCommon.Cells.Item(1571, 2).Value = "shared.cell.create.403.subtitle"
$wsCommon.Cells.Item(1571, 3).Value = "Omlouváme se, ale pro vytváření článků nemáte patřičný certifikát (nebo licenci); podívejte se na tržišti po certifikátech souvisejících se správou celého tržiště nebo správou článků."

# Row 1571 holds a long wrapped subtitle, matching row height of similar rows
$wsCommon.Rows.Item(1571).RowHeight = 26.25

# =========================================================
# Sheet: "tokens" - register the new feature token
# =========================================================
$wsTokens = $wb.Worksheets.Item("tokens")
$tokensTemplate = $wsTokens.Range("A11")
$tokensTemplate.Copy($wsTokens.Range("A12"))
$wsTokens.Cells.Item(12, 1).Value = "feature.cell.create"

# =========================================================
# Sheet: "certificates" - register the new (premium) feature
# certificates; the premium row is filled in before the base row
# =========================================================
$wsCerts = $wb.Worksheets.Item("certificates")
$certsTemplate = $wsCerts.Range("A6:D6")

$certsTemplate.Copy($wsCerts.Range("A15:D15"))
$wsCerts.Cells.Item(15, 1).Value = "feature.premium.cell.create"
$wsCerts.Cells.Item(15, 2).Value = "WHZE-SHATLV-W6SB75CL"
$wsCerts.Cells.Item(15, 3).Value = 50000
$wsCerts.Cells.Item(15, 4).Value = "feature.cell.create"

$certsTemplate.Copy($wsCerts.Range("A14:D14"))
$wsCerts.Cells.Item(14, 1).Value = "feature.cell.create"
$wsCerts.Cells.Item(14, 2).Value = "CJWP-J7E26Q-DAORQGJT"
$wsCerts.Cells.Item(14, 3).Value = ""
$wsCerts.Cells.Item(14, 4).Value = "feature.cell.create"

# =========================================================
# Back to "Translations - Common": remaining translation keys
# =========================================================
# --- rows 1572-1588 ---
$commonTemplate.Copy($wsCommon.Range("A1572:C1572"))
$wsCommon.Cells.Item(1572, 1).Value = "cs"
$wsCommon.Cells.Item(1572, 2).Value = "common.token.feature.cell.create"
$wsCommon.Cells.Item(1572, 3).Value = "Vytváření článků"

$commonTemplate.Copy($wsCommon.Range("A1573:C1573"))
$wsCommon.Cells.Item(1573, 1).Value = "cs"
$wsCommon.Cells.Item(1573, 2).Value = "certificate.feature.cell.create"
$wsCommon.Cells.Item(1573, 3).Value = "Vytváření článků"

$commonTemplate.Copy($wsCommon.Range("A1574:C1574"))
$wsCommon.Cells.Item(1574, 1).Value = "cs"
$wsCommon.Cells.Item(1574, 2).Value = "certificate.feature.premium.cell.create"
$wsCommon.Cells.Item(1574, 3).Value = "Vytváření článků - premium"

$commonTemplate.Copy($wsCommon.Range("A1575:C1575"))
$wsCommon.Cells.Item(1575, 1).Value = "cs"
$wsCommon.Cells.Item(1575, 2).Value = "common.token.certificate.feature.premium.cell.create"
$wsCommon.Cells.Item(1575, 3).Value = "Certifikát - Vytváření článků - premium"

$commonTemplate.Copy($wsCommon.Range("A1576:C1576"))
$wsCommon.Cells.Item(1576, 1).Value = "cs"
$wsCommon.Cells.Item(1576, 2).Value = "common.token.certificate.feature.cell.create"
$wsCommon.Cells.Item(1576, 3).Value = "Certifikát - Vytváření článků"

$commonTemplate.Copy($wsCommon.Range("A1577:C1577"))
$wsCommon.Cells.Item(1577, 1).Value = "cs"
$wsCommon.Cells.Item(1577, 2).Value = "shared.cell.create.name.label"
$wsCommon.Cells.Item(1577, 3).Value = "Název článku"

$commonTemplate.Copy($wsCommon.Range("A1578:C1578"))
$wsCommon.Cells.Item(1578, 1).Value = "cs"
$wsCommon.Cells.Item(1578, 2).Value = "shared.cell.create.vendorId.label"
$wsCommon.Cells.Item(1578, 3).Value = "Výrobce článku"

$commonTemplate.Copy($wsCommon.Range("A1579:C1579"))
$wsCommon.Cells.Item(1579, 1).Value = "cs"
$wsCommon.Cells.Item(1579, 2).Value = "shared.cell.create.cost.label"
$wsCommon.Cells.Item(1579, 3).Value = "Cena článku na tržišti"

$commonTemplate.Copy($wsCommon.Range("A1580:C1580"))
$wsCommon.Cells.Item(1580, 1).Value = "cs"
$wsCommon.Cells.Item(1580, 2).Value = "shared.cell.create.cost.label.tooltip"
$wsCommon.Cells.Item(1580, 3).Value = "Uveďte prosím rozumnou cenu, za kterou bude článek dostupný na tržišti."

$commonTemplate.Copy($wsCommon.Range("A1581:C1581"))
$wsCommon.Cells.Item(1581, 1).Value = "cs"
$wsCommon.Cells.Item(1581, 2).Value = "shared.cell.create.create"
$wsCommon.Cells.Item(1581, 3).Value = "Uložit"

$commonTemplate.Copy($wsCommon.Range("A1582:C1582"))
$wsCommon.Cells.Item(1582, 1).Value = "cs"
$wsCommon.Cells.Item(1582, 2).Value = "shared.cell.create.code.label"
$wsCommon.Cells.Item(1582, 3).Value = "Kód článku"

$commonTemplate.Copy($wsCommon.Range("A1583:C1583"))
$wsCommon.Cells.Item(1583, 1).Value = "cs"
$wsCommon.Cells.Item(1583, 2).Value = "shared.cell.create.voltage.label"
$wsCommon.Cells.Item(1583, 3).Value = "Pracovní napětí"

$commonTemplate.Copy($wsCommon.Range("A1584:C1584"))
$wsCommon.Cells.Item(1584, 1).Value = "cs"
$wsCommon.Cells.Item(1584, 2).Value = "shared.cell.create.voltageMax.label"
$wsCommon.Cells.Item(1584, 3).Value = "Výchozí maximální napětí"

$commonTemplate.Copy($wsCommon.Range("A1585:C1585"))
$wsCommon.Cells.Item(1585, 1).Value = "cs"
$wsCommon.Cells.Item(1585, 2).Value = "shared.cell.create.capacity.label"
$wsCommon.Cells.Item(1585, 3).Value = "Kapacita článku"

$commonTemplate.Copy($wsCommon.Range("A1586:C1586"))
$wsCommon.Cells.Item(1586, 1).Value = "cs"
$wsCommon.Cells.Item(1586, 2).Value = "shared.cell.create.drain.label"
$wsCommon.Cells.Item(1586, 3).Value = "Maximální vybíjecí proud"

$commonTemplate.Copy($wsCommon.Range("A1587:C1587"))
$wsCommon.Cells.Item(1587, 1).Value = "cs"
$wsCommon.Cells.Item(1587, 2).Value = "shared.cell.create.typeId.label"
$wsCommon.Cells.Item(1587, 3).Value = "Typ článku"

$commonTemplate.Copy($wsCommon.Range("A1588:C1588"))
$wsCommon.Cells.Item(1588, 1).Value = "cs"
$wsCommon.Cells.Item(1588, 2).Value = "shared.cell.create.success"
$wsCommon.Cells.Item(1588, 3).Value = "Článek [{{name}}] byl úspěšně vytvořen."

# =========================================================
# View state: update selections to match final editor position
# =========================================================
$wsCommon.Activate()
$wsCommon.Range("B1580").Select()

$wsTokens.Activate()
$wsTokens.Range("A7").Select()

$wsCerts.Activate()
$wsCerts.Range("B10").Select()

$wsCommon.Activate()

Write-Output "done"
